# Natmi following Dr Hou advice
# Re-compute the Ncam1 -> Robo3 ligand-receptor edge table: a third cluster
# ("ECs") and a re-labelled cluster ("M2") join the existing "FAPs"/"sCs"
# senders, expanding the 2-row result into 4 rows (A2:T5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster / Ligand symbol / Receptor symbol / Target cluster, plus the
# 16 numeric metric columns, one row per sending cluster.
$rows = @(
    @{ Row=2; Values=@(
        "ECs", "Ncam1", "Robo3", "ECs",
        2, 0.6666666666666666, 1.182981, 3.548943, 0.02832403852590813, 0.02832403852590813,
        3, 1, 2.268399666666667, 6.805199, 1, 1,
        2.683473706073, 24.151263354657, 0.02832403852590813, 0.02832403852590813
    ) },
    @{ Row=3; Values=@(
        "FAPs", "Ncam1", "Robo3", "ECs",
        3, 1, 5.178030666666667, 15.534092, 0.1239772575307637, 0.1239772575307637,
        3, 1, 2.268399666666667, 6.805199, 1, 1,
        11.74584303825645, 105.712587344308, 0.1239772575307637, 0.1239772575307637
    ) },
    @{ Row=4; Values=@(
        "M2", "Ncam1", "Robo3", "ECs",
        1, 0.3333333333333333, 0.089474, 0.268422, 0.002142270267288404, 0.002142270267288404,
        3, 1, 2.268399666666667, 6.805199, 1, 1,
        0.2029627917753333, 1.826665125978, 0.002142270267288404, 0.002142270267288404
    ) },
    @{ Row=5; Values=@(
        "sCs", "Ncam1", "Robo3", "ECs",
        3, 1, 35.31548633333333, 105.946459, 0.8455564336760397, 0.8455564336760397,
        3, 1, 2.268399666666667, 6.805199, 1, 1,
        80.10963742670455, 720.986736840341, 0.8455564336760397, 0.8455564336760397
    ) }
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($entry in $rows) {
    $r = $entry.Row
    $vals = $entry.Values
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$r").Value = $vals[$i]
    }
}
